# "new sample file added for goods and sample"
# Add four new trailing columns (Weight, Length, Breadth, Height) to the
# product/goods sample header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = "Weight"
$ws.Range("N1").Value = "Length"
$ws.Range("O1").Value = "Breadth"
$ws.Range("P1").Value = "Height"

$ws.Range("P2").Select()
